# Insert a new weekly price record for "Espárragos" (Macroferia Regional de
# Talca) as row 107, pushing the existing rows 107-121 down to 108-122.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(107).Insert()

$ws.Range("A107").Value = 5
$ws.Range("B107").Value = "Macroferia Regional de Talca"
$ws.Range("C107").Value = "Maule"
$ws.Range("D107").Value = 45223
$ws.Range("E107").Value = 7
$ws.Range("F107").Value = 300000000
$ws.Range("G107").Value = "Espárragos"
$ws.Range("H107").Value = "Sin especificar"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1100
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1140
$ws.Range("N107").Value = "`$/kilo"
$ws.Range("O107").Value = "Provincia de Linares"
$ws.Range("P107").Value = 1140
$ws.Range("Q107").Value = 1
$ws.Range("R107").Value = "Hortaliza"
